$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the target cells to Text format first so numeric-looking strings
# (e.g. "1.005") are stored as text, matching the source workbook, then
# clear the temporary formatting so no style index is left behind.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "27.403.56"
$ws.Range("E2").Value = "  -1.63%  "
$ws.Range("D3").Value = "1.711.65"
$ws.Range("E3").Value = "  -1.72%  "
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "224.46"
$ws.Range("E5").Value = "  -1.44%  "
$ws.Range("D6").Value = "0.5339"
$ws.Range("E6").Value = "  -2.25%  "
$ws.Range("D7").Value = "1.005"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "0.2668"
$ws.Range("E8").Value = "  -3.81%  "
$ws.Range("D9").Value = "0.06606"
$ws.Range("E9").Value = "  -1.98%  "
$ws.Range("E10").Value = "  -4.38%  "
$ws.Range("D11").Value = "0.07651"
$ws.Range("E11").Value = "  -1.65%  "
$ws.Range("D12").Value = "4.575"
$ws.Range("E12").Value = "  -2.64%  "
$ws.Range("D13").Value = "1.745.97"
$ws.Range("E13").Value = "  +0.04%  "
$ws.Range("D14").Value = "1.948.70"
$ws.Range("E14").Value = "  -1.65%  "
$ws.Range("D15").Value = "0.5769"
$ws.Range("E15").Value = "  -3.46%  "
$ws.Range("D16").Value = "0.0₅8177"
$ws.Range("E16").Value = "  -2.62%  "
$ws.Range("D17").Value = "67.96"
$ws.Range("E17").Value = "  -1.45%  "
$ws.Range("D18").Value = "27.383.96"
$ws.Range("E18").Value = "  -1.69%  "
$ws.Range("D19").Value = "216.17"
$ws.Range("E19").Value = "  -3.64%  "
$ws.Range("E20").Value = "  +0.12%  "
$ws.Range("D21").Value = "4.678"
$ws.Range("E21").Value = "  -3.42%  "
$ws.Range("E22").Value = "  -4.12%  "
$ws.Range("D23").Value = "5.984"
$ws.Range("E23").Value = "  -4.15%  "
$ws.Range("D24").Value = "1.005"
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("D25").Value = "142.28"
$ws.Range("E25").Value = "  -2.75%  "
$ws.Range("D26").Value = "1.726"
$ws.Range("E26").Value = "  +2.58%  "
$ws.Range("D27").Value = "0.1219"
$ws.Range("E27").Value = "  -2.64%  "
$ws.Range("D28").Value = "7.285"
$ws.Range("E28").Value = "  -2.37%  "
$ws.Range("D29").Value = "16.35"
$ws.Range("E29").Value = "  -4.94%  "
$ws.Range("D30").Value = "0.05418"
$ws.Range("E30").Value = "  -4.51%  "
$ws.Range("D31").Value = "1.295"
$ws.Range("E31").Value = "  -1.27%  "
$ws.Range("D32").Value = "3.513"
$ws.Range("E32").Value = "  -5.03%  "
$ws.Range("D33").Value = "3.438"
$ws.Range("E33").Value = "  -2.41%  "
$ws.Range("E34").Value = "  -2.22%  "
$ws.Range("D35").Value = "2.882"
$ws.Range("E35").Value = "  +0.89%  "
$ws.Range("D36").Value = "0.9503"
$ws.Range("E36").Value = "  -2.80%  "
$ws.Range("D37").Value = "2.420"
$ws.Range("E37").Value = "  -1.25%  "
$ws.Range("D38").Value = "0.5874"
$ws.Range("E38").Value = "  -1.69%  "
$ws.Range("D39").Value = "0.01634"
$ws.Range("E39").Value = "  -1.93%  "
$ws.Range("D40").Value = "5.877"
$ws.Range("E40").Value = "  -1.90%  "
$ws.Range("D41").Value = "1.044.27"
$ws.Range("E41").Value = "  -0.23%  "
$ws.Range("D42").Value = "1.005"
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("D43").Value = "0.8426"
$ws.Range("E43").Value = "  -1.06%  "
$ws.Range("D44").Value = "101.02"
$ws.Range("E44").Value = "  -0.98%  "
$ws.Range("D45").Value = "1.854.36"
$ws.Range("E45").Value = "  -1.75%  "
$ws.Range("E46").Value = "  +7.65%  "
$ws.Range("D47").Value = "58.15"
$ws.Range("E47").Value = "  -2.46%  "
$ws.Range("D48").Value = "0.4514"
$ws.Range("E48").Value = "  +1.71%  "
$ws.Range("D49").Value = "1.005"
$ws.Range("E49").Value = "  +0.24%  "
$ws.Range("D50").Value = "8.095"
$ws.Range("E50").Value = "  -2.33%  "
$ws.Range("D51").Value = "0.05242"
$ws.Range("E51").Value = "  -1.60%  "

$ws.Range("D2:E51").ClearFormats()
